$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 22080
$ws.Range("E2").Value = 1328
$ws.Range("F2").Value = 1328
$ws.Range("G2").Value = -791
$ws.Range("H2").Value = -686
$ws.Range("I2").Value = -686
$ws.Range("K2").Value = 51331
$ws.Range("L2").Value = 31705
$ws.Range("M2").Value = 19626
$ws.Range("N2").Value = 19626
$ws.Range("P2").Value = 3996
$ws.Range("Q2").Value = -49
$ws.Range("R2").Value = -1668
$ws.Range("S2").Value = 1810
$ws.Range("T2").Value = 266
$ws.Range("U2").Value = -315
$ws.Range("V2").Value = 16781
$ws.Range("W2").Value = 6.02
$ws.Range("X2").Value = -3.11
$ws.Range("Y2").Value = -3.45
$ws.Range("Z2").Value = -1.36
$ws.Range("AA2").Value = 161.55
$ws.Range("AB2").Value = 351.68
$ws.Range("AC2").Value = -842
$ws.Range("AD2").Value = -11.17
$ws.Range("AE2").Value = 24105
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 94
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = -80.36
$ws.Range("AJ2").Value = 58693231

# Row 3
$ws.Range("D3").Value = 11853
$ws.Range("E3").Value = -1279
$ws.Range("F3").Value = -1279
$ws.Range("G3").Value = -4313
$ws.Range("H3").Value = -5207
$ws.Range("I3").Value = -5207
$ws.Range("K3").Value = 42257
$ws.Range("L3").Value = 28114
$ws.Range("M3").Value = 14143
$ws.Range("N3").Value = 14143
$ws.Range("P3").Value = 4206
$ws.Range("Q3").Value = 2150
$ws.Range("R3").Value = 977
$ws.Range("S3").Value = -3722
$ws.Range("T3").Value = 175
$ws.Range("U3").Value = 1976
$ws.Range("V3").Value = 13802
$ws.Range("W3").Value = -10.79
$ws.Range("X3").Value = -43.93
$ws.Range("Y3").Value = -30.84
$ws.Range("Z3").Value = -11.13
$ws.Range("AA3").Value = 198.78
$ws.Range("AB3").Value = 205.56
$ws.Range("AC3").Value = -6189
$ws.Range("AD3").Value = -0.78
$ws.Range("AE3").Value = 16710
$ws.Range("AF3").Value = 0.29
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = -4.99
$ws.Range("AJ3").Value = 63160470

# Row 4
$ws.Range("D4").Value = 13342
$ws.Range("E4").Value = 198
$ws.Range("F4").Value = 198
$ws.Range("G4").Value = -2367
$ws.Range("H4").Value = -3570
$ws.Range("I4").Value = -3570
$ws.Range("K4").Value = 30300
$ws.Range("L4").Value = 19654
$ws.Range("M4").Value = 10647
$ws.Range("N4").Value = 10647
$ws.Range("P4").Value = 543
$ws.Range("Q4").Value = -1005
$ws.Range("R4").Value = 5374
$ws.Range("S4").Value = -4653
$ws.Range("T4").Value = 64
$ws.Range("U4").Value = -1069
$ws.Range("V4").Value = 9068
$ws.Range("W4").Value = 1.49
$ws.Range("X4").Value = -26.76
$ws.Range("Y4").Value = -28.8
$ws.Range("Z4").Value = -9.84
$ws.Range("AA4").Value = 184.6
$ws.Range("AB4").Value = 1649.98
$ws.Range("AC4").Value = -4039
$ws.Range("AD4").Value = -0.86
$ws.Range("AE4").Value = 11943
$ws.Range("AF4").Value = 0.29
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 64997543

# Row 5
$ws.Range("D5").Value = 15359
$ws.Range("E5").Value = 589
$ws.Range("F5").Value = 589
$ws.Range("G5").Value = -1348
$ws.Range("H5").Value = -1840
$ws.Range("I5").Value = -1845
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 28456
$ws.Range("L5").Value = 18802
$ws.Range("M5").Value = 9654
$ws.Range("N5").Value = 8850
$ws.Range("O5").Value = 804
$ws.Range("P5").Value = 545
$ws.Range("Q5").Value = 553
$ws.Range("R5").Value = -423
$ws.Range("S5").Value = 183
$ws.Range("T5").Value = 58
$ws.Range("U5").Value = 495
$ws.Range("V5").Value = 8635
$ws.Range("W5").Value = 3.84
$ws.Range("X5").Value = -11.98
$ws.Range("Y5").Value = -18.93
$ws.Range("Z5").Value = -6.26
$ws.Range("AA5").Value = 194.75
$ws.Range("AB5").Value = 1361.24
$ws.Range("AC5").Value = -1990
$ws.Range("AD5").Value = -1.39
$ws.Range("AE5").Value = 9872
$ws.Range("AF5").Value = 0.28
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 65499044

# Row 6
$ws.Range("D6").Value = 15478
$ws.Range("E6").Value = -522
$ws.Range("F6").Value = -522
$ws.Range("G6").Value = -5291
$ws.Range("H6").Value = -5518
$ws.Range("I6").Value = -5527
$ws.Range("K6").Value = 23991
$ws.Range("L6").Value = 20314
$ws.Range("M6").Value = 3677
$ws.Range("N6").Value = 2870
$ws.Range("P6").Value = 604
$ws.Range("Q6").Value = 945
$ws.Range("R6").Value = -1203
$ws.Range("S6").Value = 139
$ws.Range("T6").Value = 43
$ws.Range("U6").Value = 902
$ws.Range("V6").Value = 8615
$ws.Range("W6").Value = -3.37
$ws.Range("X6").Value = -35.65
$ws.Range("Y6").Value = -94.32
$ws.Range("Z6").Value = -21.04
$ws.Range("AA6").Value = 552.5
$ws.Range("AB6").Value = 268.37
$ws.Range("AC6").Value = -5484
$ws.Range("AD6").Value = -0.24
$ws.Range("AE6").Value = 2810
$ws.Range("AF6").Value = 0.47
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 77987263

# Structural removals (cells deleted entirely)
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()

Write-Host "edit complete"